$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert the new columns.
#    Original layout:  A symbol | B date | C is_today | D close | E today_support
#                       | F support | G message_s | H message_e0 | I message_e1
#                       | J message_e2 | K rsi | L cci | M sma5 | N sma50 | O sma200
#                       | P ema5 | Q ema20 | R..AB (flags)
#    New layout adds:  is_first_buy_yn before close, csp_bullish_candle +
#                       volume_inconsistency_alert before message_s, and
#                       ema5_flag appended as the new last column.
# ---------------------------------------------------------------------------
$ws.Range("D1").EntireColumn.Insert()
$ws.Range("H1:I1").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 2. Fill in the header text for the newly inserted / appended columns.
# ---------------------------------------------------------------------------
$ws.Range("D1").Value2 = "is_first_buy_yn"
$ws.Range("H1").Value2 = "csp_bullish_candle"
$ws.Range("I1").Value2 = "volume_inconsistency_alert"
$ws.Range("AF1").Value2 = "ema5_flag"

# ---------------------------------------------------------------------------
# 3. Formatting.
# ---------------------------------------------------------------------------
# D1 and AF1 pick up the plain header style (same as A1/B1/C1 etc.).
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("AF1").PasteSpecial(-4122)

# E1 ("close") now also gets the orange highlight fill used on ema5/ema20,
# reusing the already-existing currency+orange style (S1/T1 after the shift).
$ws.Range("S1").Copy()
$ws.Range("E1").PasteSpecial(-4122)

# H1 / I1 (new) get a plain wrapped header style first (clears the currency
# format the column Insert inherited from its left neighbour), then a
# yellow highlight fill on top.
$ws.Range("A1").Copy()
$ws.Range("H1:I1").PasteSpecial(-4122)
$ws.Range("H1:I1").Interior.Color = 65535

# J1 ("message_s") gets the same treatment with the orange highlight fill.
$ws.Range("A1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Interior.Color = 49407

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. AutoFilter now covers the wider range.
# ---------------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:AF1").AutoFilter()

# ---------------------------------------------------------------------------
# 5. Conditional formatting: the existing color-scale rule now targets
#    U1:AE1048576 (shifted), and a brand-new color-scale rule is added for
#    the newly appended AF column, placed at the top priority.
# ---------------------------------------------------------------------------
$oldRuleRange = $ws.Range("R1:AB1048576")
$oldRule = $oldRuleRange.FormatConditions.Item(1)
$oldRule.ModifyAppliesToRange($ws.Range("U1:AE1048576"))

$newRuleRange = $ws.Range("AF1:AF1048576")
$newRule = $newRuleRange.FormatConditions.AddColorScale(3)
$newRule.SetFirstPriority()

# ---------------------------------------------------------------------------
# 6. Defined name (_FilterDatabase) and zoom level.
# ---------------------------------------------------------------------------
$wb.Names.Item("watchlist!_FilterDatabase").RefersTo = "=watchlist!`$A`$1:`$AF`$1"

$excel.ActiveWindow.Zoom = 55
